$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Förändrad" (Changed) date column C for all existing data rows (2..358) ---
for ($r = 2; $r -le 358; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# --- Give row 358 the explicit 15pt custom row height (as in the target) ---
$ws.Rows.Item(358).RowHeight = 15

# --- Append four new data rows (359..362) ---
function Add-Row($row, $beteckning, $datum, $forandrad, $area) {
    $ws.Cells.Item($row, 1).Value = $beteckning
    $ws.Cells.Item($row, 2).Value = $datum
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 3).Value = $forandrad
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = "BLEKINGE LÄN"
    $ws.Cells.Item($row, 5).Value = "OLOFSTRÖM"
    $ws.Cells.Item($row, 7).Value = $area
    for ($c = 8; $c -le 17; $c++) {
        $ws.Cells.Item($row, $c).Value = 0
    }
    # Column R stays an empty, wrap-text styled cell like the rest of the table
    $ws.Cells.Item(358, 18).Copy($ws.Cells.Item($row, 18))
}

Add-Row 359 "A 45688-2023" 45195 45202 6.7
Add-Row 360 "A 45689-2023" 45195 45202 3.4
Add-Row 361 "A 46531-2023" 45197 45202 2.7
Add-Row 362 "A 46361-2023" 45197 45202 3.5

# Rows 359-361 keep the explicit 15pt row height; row 362 (the new last row) uses default height
$ws.Rows.Item(359).RowHeight = 15
$ws.Rows.Item(360).RowHeight = 15
$ws.Rows.Item(361).RowHeight = 15
